$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 1.73
$ws.Cells.Item(2, 8).Value = 5.8
$ws.Cells.Item(2, 9).Value = 6.6
$ws.Cells.Item(2, 10).Value = 3.85
$ws.Cells.Item(2, 11).Value = 4.3
$ws.Cells.Item(2, 12).Value = 1.44
$ws.Cells.Item(2, 14).Value = 3.85
$ws.Cells.Item(2, 15).Value = 1.31
$ws.Cells.Item(2, 16).Value = 1.96
$ws.Cells.Item(2, 17).Value = 1.98
$ws.Cells.Item(2, 18).Value = 1.36
$ws.Cells.Item(2, 19).Value = 3.5
$ws.Cells.Item(2, 20).Value = 1.93
$ws.Cells.Item(2, 21).Value = 1.97
$ws.Cells.Item(2, 23).Value = 2.4
$ws.Cells.Item(2, 24).Value = 19.5
$ws.Cells.Item(2, 25).Value = 970
$ws.Cells.Item(2, 28).Value = 9.199999999999999
$ws.Cells.Item(2, 29).Value = 11
$ws.Cells.Item(2, 30).Value = 55
$ws.Cells.Item(2, 32).Value = 12.5
$ws.Cells.Item(2, 33).Value = 9
$ws.Cells.Item(2, 34).Value = 22
$ws.Cells.Item(2, 36).Value = 180
$ws.Cells.Item(2, 37).Value = 19
$ws.Cells.Item(2, 38).Value = 150
$ws.Cells.Item(2, 40).Value = 10.5
# Row 3
$ws.Cells.Item(3, 6).Value = 1.7
$ws.Cells.Item(3, 7).Value = 1.87
$ws.Cells.Item(3, 8).Value = 4.5
$ws.Cells.Item(3, 9).Value = 5.8
$ws.Cells.Item(3, 10).Value = 3.9
$ws.Cells.Item(3, 11).Value = 4.7
$ws.Cells.Item(3, 12).Value = 1.32
$ws.Cells.Item(3, 14).Value = 4.5
$ws.Cells.Item(3, 16).Value = 2.22
$ws.Cells.Item(3, 17).Value = 1.71
$ws.Cells.Item(3, 18).Value = 1.48
$ws.Cells.Item(3, 19).Value = 2.66
$ws.Cells.Item(3, 21).Value = 2.18
$ws.Cells.Item(3, 23).Value = 2.16
$ws.Cells.Item(3, 24).Value = 22
$ws.Cells.Item(3, 26).Value = 95
$ws.Cells.Item(3, 31).Value = 220
$ws.Cells.Item(3, 32).Value = 12.5
$ws.Cells.Item(3, 34).Value = 18.5
$ws.Cells.Item(3, 37).Value = 18
$ws.Cells.Item(3, 39).Value = 300
$ws.Cells.Item(3, 40).Value = 9.4
$ws.Cells.Item(3, 41).Value = 110
# Row 4
$ws.Cells.Item(4, 13).Value = 1.04
$ws.Cells.Item(4, 20).Value = 1.69
# Row 5
$ws.Cells.Item(5, 7).Value = 1.4
$ws.Cells.Item(5, 8).Value = 10
$ws.Cells.Item(5, 10).Value = 5.4
$ws.Cells.Item(5, 11).Value = 5.6
$ws.Cells.Item(5, 12).Value = 1.38
$ws.Cells.Item(5, 13).Value = 1.05
$ws.Cells.Item(5, 14).Value = 4.2
$ws.Cells.Item(5, 16).Value = 2.12
$ws.Cells.Item(5, 18).Value = 1.42
$ws.Cells.Item(5, 19).Value = 3.25
$ws.Cells.Item(5, 20).Value = 2.18
$ws.Cells.Item(5, 21).Value = 1.77
$ws.Cells.Item(5, 22).Value = 1.1
$ws.Cells.Item(5, 23).Value = 3.5
$ws.Cells.Item(5, 24).Value = 18.5
$ws.Cells.Item(5, 25).Value = 32
$ws.Cells.Item(5, 40).Value = 6.8
# Row 6
$ws.Cells.Item(6, 6).Value = 4.5
$ws.Cells.Item(6, 9).Value = 2.08
$ws.Cells.Item(6, 12).Value = 1.53
$ws.Cells.Item(6, 19).Value = 4.8
$ws.Cells.Item(6, 20).Value = 2.08
$ws.Cells.Item(6, 22).Value = 1.92
$ws.Cells.Item(6, 23).Value = 1.26
$ws.Cells.Item(6, 27).Value = 25
$ws.Cells.Item(6, 29).Value = 7.6
$ws.Cells.Item(6, 31).Value = 25
$ws.Cells.Item(6, 34).Value = 22
$ws.Cells.Item(6, 40).Value = 95
$ws.Cells.Item(6, 41).Value = 22
# Row 7
$ws.Cells.Item(7, 6).Value = 2.28
$ws.Cells.Item(7, 7).Value = 2.34
$ws.Cells.Item(7, 8).Value = 3.6
$ws.Cells.Item(7, 10).Value = 3.35
$ws.Cells.Item(7, 13).Value = 1.09
$ws.Cells.Item(7, 14).Value = 3.3
$ws.Cells.Item(7, 16).Value = 1.77
$ws.Cells.Item(7, 17).Value = 2.2
$ws.Cells.Item(7, 18).Value = 1.28
$ws.Cells.Item(7, 19).Value = 4.2
$ws.Cells.Item(7, 20).Value = 1.86
$ws.Cells.Item(7, 21).Value = 2.02
$ws.Cells.Item(7, 22).Value = 1.36
$ws.Cells.Item(7, 23).Value = 1.74
$ws.Cells.Item(7, 30).Value = 15
$ws.Cells.Item(7, 39).Value = 140
# Row 8
$ws.Cells.Item(8, 6).Value = 2.52
$ws.Cells.Item(8, 12).Value = 1.53
$ws.Cells.Item(8, 16).Value = 1.67
$ws.Cells.Item(8, 17).Value = 2.44
# Row 9
$ws.Cells.Item(9, 6).Value = 2.68
$ws.Cells.Item(9, 7).Value = 2.82
$ws.Cells.Item(9, 8).Value = 2.98
$ws.Cells.Item(9, 9).Value = 3.15
$ws.Cells.Item(9, 12).Value = 1.55
$ws.Cells.Item(9, 13).Value = 1.11
$ws.Cells.Item(9, 14).Value = 2.86
$ws.Cells.Item(9, 15).Value = 1.5
$ws.Cells.Item(9, 16).Value = 1.6
$ws.Cells.Item(9, 17).Value = 2.52
$ws.Cells.Item(9, 22).Value = 1.47
$ws.Cells.Item(9, 23).Value = 1.54
$ws.Cells.Item(9, 25).Value = 9.4
$ws.Cells.Item(9, 26).Value = 22
$ws.Cells.Item(9, 28).Value = 8.800000000000001
$ws.Cells.Item(9, 29).Value = 7.4
$ws.Cells.Item(9, 30).Value = 14
$ws.Cells.Item(9, 31).Value = 1000
$ws.Cells.Item(9, 32).Value = 20
$ws.Cells.Item(9, 33).Value = 13.5
$ws.Cells.Item(9, 37).Value = 980
$ws.Cells.Item(9, 39).Value = 180
# Row 10
$ws.Cells.Item(10, 6).Value = 2.4
$ws.Cells.Item(10, 7).Value = 2.48
$ws.Cells.Item(10, 8).Value = 3.6
$ws.Cells.Item(10, 9).Value = 3.75
$ws.Cells.Item(10, 12).Value = 1.62
$ws.Cells.Item(10, 15).Value = 1.59
$ws.Cells.Item(10, 16).Value = 1.52
$ws.Cells.Item(10, 17).Value = 2.8
$ws.Cells.Item(10, 18).Value = 1.18
$ws.Cells.Item(10, 19).Value = 6
$ws.Cells.Item(10, 21).Value = 1.75
$ws.Cells.Item(10, 22).Value = 1.36
$ws.Cells.Item(10, 23).Value = 1.67
$ws.Cells.Item(10, 25).Value = 9.6
$ws.Cells.Item(10, 26).Value = 27
$ws.Cells.Item(10, 27).Value = 100
$ws.Cells.Item(10, 30).Value = 21
$ws.Cells.Item(10, 31).Value = 1000
$ws.Cells.Item(10, 32).Value = 16
$ws.Cells.Item(10, 34).Value = 32
$ws.Cells.Item(10, 36).Value = 42
$ws.Cells.Item(10, 40).Value = 1000
$ws.Cells.Item(10, 41).Value = 110
# Row 11
$ws.Cells.Item(11, 12).Value = 1.48
$ws.Cells.Item(11, 14).Value = 3.3
$ws.Cells.Item(11, 16).Value = 1.77
$ws.Cells.Item(11, 17).Value = 2.22
$ws.Cells.Item(11, 20).Value = 1.89
$ws.Cells.Item(11, 22).Value = 1.4
$ws.Cells.Item(11, 34).Value = 19.5
$ws.Cells.Item(11, 41).Value = 60
# Row 12
$ws.Cells.Item(12, 7).Value = 3.35
$ws.Cells.Item(12, 8).Value = 2.6
$ws.Cells.Item(12, 9).Value = 2.64
$ws.Cells.Item(12, 12).Value = 1.58
$ws.Cells.Item(12, 14).Value = 2.74
$ws.Cells.Item(12, 17).Value = 2.66
$ws.Cells.Item(12, 19).Value = 5.5
$ws.Cells.Item(12, 20).Value = 2.14
$ws.Cells.Item(12, 21).Value = 1.8
$ws.Cells.Item(12, 22).Value = 1.6
$ws.Cells.Item(12, 29).Value = 7
$ws.Cells.Item(12, 36).Value = 65
$ws.Cells.Item(12, 38).Value = 80
$ws.Cells.Item(12, 40).Value = 1000
# Row 13
$ws.Cells.Item(13, 6).Value = 2.1
$ws.Cells.Item(13, 8).Value = 3.3
$ws.Cells.Item(13, 11).Value = 4.2
$ws.Cells.Item(13, 12).Value = 1.41
$ws.Cells.Item(13, 14).Value = 3.35
$ws.Cells.Item(13, 16).Value = 1.81
$ws.Cells.Item(13, 17).Value = 1.96
$ws.Cells.Item(13, 19).Value = 3.45
$ws.Cells.Item(13, 20).Value = 1.73
$ws.Cells.Item(13, 21).Value = 2.02
$ws.Cells.Item(13, 22).Value = 1.32
$ws.Cells.Item(13, 29).Value = 1000
# Row 14
$ws.Cells.Item(14, 6).Value = 2.18
$ws.Cells.Item(14, 7).Value = 2.3
$ws.Cells.Item(14, 8).Value = 3.45
$ws.Cells.Item(14, 9).Value = 3.75
$ws.Cells.Item(14, 10).Value = 3.5
$ws.Cells.Item(14, 11).Value = 3.8
$ws.Cells.Item(14, 12).Value = 1.49
$ws.Cells.Item(14, 13).Value = 1.08
$ws.Cells.Item(14, 14).Value = 3.15
$ws.Cells.Item(14, 15).Value = 1.42
$ws.Cells.Item(14, 16).Value = 1.76
$ws.Cells.Item(14, 17).Value = 2.26
$ws.Cells.Item(14, 18).Value = 1.27
$ws.Cells.Item(14, 19).Value = 4.3
$ws.Cells.Item(14, 20).Value = 1.84
$ws.Cells.Item(14, 21).Value = 1.94
$ws.Cells.Item(14, 22).Value = 1.37
$ws.Cells.Item(14, 23).Value = 1.76
$ws.Cells.Item(14, 24).Value = 15
$ws.Cells.Item(14, 25).Value = 14.5
$ws.Cells.Item(14, 28).Value = 10.5
$ws.Cells.Item(14, 29).Value = 9.800000000000001
$ws.Cells.Item(14, 30).Value = 17.5
$ws.Cells.Item(14, 31).Value = 60
$ws.Cells.Item(14, 33).Value = 11.5
$ws.Cells.Item(14, 34).Value = 22
$ws.Cells.Item(14, 36).Value = 32
$ws.Cells.Item(14, 37).Value = 32
$ws.Cells.Item(14, 38).Value = 55
$ws.Cells.Item(14, 40).Value = 28
$ws.Cells.Item(14, 41).Value = 1000
